# 19/12/2025: Update the list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-15 (Job ID, Company, Job Title, Candidate, Status)
$data = @(
    @(580, "Legion Security", "Sales Engineer (US)", "Michael Saladino", "2nd Interview"),
    @(663, "Blink Ops", "Sales Engineer UK", "Adam Evans", "4th Interview"),
    @(663, "Blink Ops", "Sales Engineer UK", "Graham Rance", "4th Interview"),
    @(663, "Blink Ops", "Sales Engineer UK", "Kev Pyart", "4th Interview"),
    @(663, "Blink Ops", "Sales Engineer UK", "LLOYD WEBB", "1st Interview"),
    @(673, "Redwood Software", "SE UK", "Joseph Falvey", "1st Interview"),
    @(673, "Redwood Software", "SE UK", "Navid Ghavami", "1st Interview"),
    @(673, "Redwood Software", "SE UK", "RICHARD JUDD", "3rd Interview"),
    @(707, "Dash0", "Sales Engineer EMEA (UK, Nordics, Benelux, Germany) x 2", "Harry Kimpel", "2nd Interview"),
    @(707, "Dash0", "Sales Engineer EMEA (UK, Nordics, Benelux, Germany) x 2", "Patrick Schrimpf", "3rd Interview"),
    @(730, "PointFive", "PointFive SE EST", "Matthew Hughes", "1st Interview"),
    @(768, "Adaptive6", "Senior Sales Engineer (US)", "Matthew Hughes", "2nd Interview"),
    @(849, "Oasis Security", "Director of Product Marketing", "Alex Spiliotes", "1st Interview"),
    @(866, "CyCognito", "BDR US", "Bryce Morais", "CV Sent")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}

# Remove the now-obsolete trailing rows (previously rows 16-19)
$ws.Range("A16:E19").Delete()
